$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix wording in two existing 'Negative Aspects' cells (remove the leading
#     'Once more, ' / 'Again, ' filler phrase) ---
$ws.Range("E4").Value = 'To examine cohort demographics, downloading the database is necessary. Additionally, it solely covers ICU data, lacking information on medical devices.'
$ws.Range("E6").Value = 'To analyze cohort demographics, downloading the database is required. Furthermore, it exclusively encompasses ICU data, with no information provided on medical devices.'

# --- Append three new database rows (10-12) with their reference material ---
# Row 10
$ws.Range("A10").Value = 'Data visualization tool MIMIC'
$ws.Range("B10").Value = 'Lee, J., Ribey, E., \& Wallace, J. R. (2015). A web-based data visualization tool for the MIMIC-II database. BMC medical informatics and decision making, 16, 1-8.'
$ws.Range("C10").Value = 'The tool provides two primary functions: Explore and Compare. With the Explore feature, users can choose a patient group within MIMIC-II and view the distributions of different administrative, demographic, and clinical variables within that group. The Compare feature allows users to select two patient groups and visually analyze the differences between them across various variables. This tool proves valuable to seasoned MIMIC-II researchers by streamlining the laborious process of writing SQL queries and manually visualizing extracted data, thereby significantly expediting their research efforts.'
$ws.Range("D10").Value = 'The tool features an Explore option, which enables users to select a patient cohort based on factors such as admission ICU service type, gender, age, and diagnosis. This functionality aligns with the objectives of the catalog. Additionally, it allows for visualization of the selected cohorts and potentially facilitates comparison. Thus, it offers a preview of the quantity of available data.'
$ws.Range("E10").Value = 'However, it does not allow for visualizing the data journey or the data source.'
$ws.Rows.Item(10).RowHeight = 186.75

# Row 11
$ws.Range("A11").Value = 'Interactive data visualization tool to assess the impact of decision support on clinical operations'
$ws.Range("B11").Value = 'Huber, T. C., Krishnaraj, A., ... , \& Gaskin, C. M. (2018). Developing an interactive data visualization tool to assess the impact of decision support on clinical operations. Journal of digital imaging, 31, 640-645.'
$ws.Range("C11").Value = 'Radiology practices nationwide are adopting clinical decision support (CDS) software. This software aids referring providers with imaging study decisions at the point of order entry, generating a significant volume of data ripe for research and quality improvement. To better understand and analyze trends within this data, an interactive data visualization dashboard was developed using a commercially available platform. By integrating a clinical decision support product into the electronic health record and exporting its data to the visualization platform (Tableau, Seattle, WA), real-time visualization of CDS-generated data became possible. This dashboard enhances the analysis of CDS output, facilitating hypothesis generation and streamlining research and quality improvement endeavors. Integrating data visualization tools with clinical decision support systems simplifies data analysis and enhances the efficiency of research and quality improvement processes.'
$ws.Range("D11").Value = 'This article describes a tool very similar to what we would like to develop, but without the online aspect. It give some details on the journey of the data and it enables high-level visualization of a set of parameters relevant to the tool''s objective.'
$ws.Range("E11").Value = 'This tool lacks filtering capabilities, despite having elements to modify data display on the dashboard. However, it doesn''t address our specific question due to its divergent objectives.'
$ws.Rows.Item(11).RowHeight = 309.75

# Row 12
$ws.Range("A12").Value = 'Data visualization tool Ophtalmology'
$ws.Range("B12").Value = 'Kortüm, K. U., Müller, M., ... \& Hirneiss, C. (2017). Using electronic health records to build an ophthalmologic data warehouse and visualize patients'' data. American journal of ophthalmology, 178, 84-93.'
$ws.Range("C12").Value = 'An academic ophthalmologic center established a near-real-time data warehouse (DW) to leverage the growing digital data from electronic medical records (EMR) and diagnostic devices. They integrated specific macular clinic interfaces into the hospital information system, enabling seamless ordering of imaging modalities. The DW, powered by an SQL database, compiled data from over 325,767 patients since 2002 and included a data discovery tool. Notably, a search for patients with age-related macular degeneration who underwent cataract surgery and received at least 10 intravitreal injections yielded 450 patients meeting the criteria. A web-based browsing tool facilitated data visualization and filtering based on various criteria, streamlining analysis and enhancing insights into ophthalmologic data.'
$ws.Range("D12").Value = 'This tool aligns with our vision for developing our catalog from the clinical data center''s database. It enables targeting patient cohorts, diagnoses, and procedures. Additionally, a dashboard-style visualization appears to be available. It seems that the origin of the images (imagers) is also provided.'
$ws.Range("E12").Value = 'This tool facilitates the extraction of a significant amount of data but is limited to the field of ophthalmology. Moreover, its objectives differ as it aims to simplify the identification of eligible patients for studies.'
$ws.Rows.Item(12).RowHeight = 258

# --- Column width adjustments (A widened to fit new long header, C widened too) ---
$ws.Columns.Item(1).ColumnWidth = 24.666666666666668
$ws.Columns.Item(3).ColumnWidth = 46

# --- Move selection to the last filled cell and scroll the new rows into view ---
$ws.Range("E12").Select()
